$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Fix typo: "visibel" -> "visible" (row 8, column D) - edited in place first
$ws.Cells.Item(8, 4).Value = "Make browse playlist mode visible"

# New row 9 - Issue #8 "Add playlist"
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "DONE"
$ws.Cells.Item(9, 4).Value = "Add playlist"
$ws.Cells.Item(9, 5).Value = "Tidy UI"
$ws.Cells.Item(9, 6).Value = "make heading h4 and buttons same as on home"

# New row 10 - Issue #9 "Cant save newly created playlist"
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 3).Value = "BUG"

# Add Type column entries for existing rows 5-8 and the new row 9
$ws.Cells.Item(5, 3).Value = "UI"
$ws.Cells.Item(6, 3).Value = "Arch"
$ws.Cells.Item(7, 3).Value = "UI"
$ws.Cells.Item(8, 3).Value = "UI"
$ws.Cells.Item(9, 3).Value = "UI"

# Finish row 10
$ws.Cells.Item(10, 6).Value = "Cant save newly created playlist"

# Row 6 no longer needs the taller (two-line) height - shrink back to default
$ws.Rows.Item(6).EntireRow.AutoFit()

# Update the active selection to match the edited workbook
$ws.Range("D9").Select()
